$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns for this run's cryptos snapshot.
# Cell values are text (e.g. "28.01", "0.999", "  +0.11%  "), not numbers, so we
# temporarily force a text number format while assigning, then restore the default
# "Normal" style so we don't leave a stray text-format style on the cell.
$cellUpdates = @{
    "D2" = "68.429.83"
    "D3" = "2.649.32"
    "E3" = "  +0.11%  "
    "E4" = "  -0.01%  "
    "D5" = "597.64"
    "E5" = "  -0.06%  "
    "D6" = "159.37"
    "E6" = "  +2.94%  "
    "E7" = "  +0.02%  "
    "E8" = "  -1.33%  "
    "D9" = "2.647.98"
    "E9" = "  +0.08%  "
    "E10" = "  -1.36%  "
    "E11" = "  -1.11%  "
    "E12" = "  +0.36%  "
    "E13" = "  -1.25%  "
    "D14" = "28.01"
    "E14" = "  -0.11%  "
    "D15" = "3.131.25"
    "E15" = "  +0.13%  "
    "D16" = "0.0000188"
    "E16" = "  -2.82%  "
    "D17" = "68.292.06"
    "E17" = "  -0.02%  "
    "D18" = "2.669.88"
    "E18" = "  +0.89%  "
    "D19" = "11.48"
    "E19" = "  +1.00%  "
    "D20" = "363.75"
    "E20" = "  -0.19%  "
    "D21" = "7.45"
    "E21" = "  -0.53%  "
    "D22" = "4.42"
    "E22" = "  +0.82%  "
    "D23" = "4.77"
    "E23" = "  -2.55%  "
    "D24" = "2.08"
    "E24" = "  +0.44%  "
    "E25" = "  -0.47%  "
    "E26" = "  +0.02%  "
    "D27" = "9.85"
    "E27" = "  +0.40%  "
    "D28" = "2.780.48"
    "E28" = "  +0.15%  "
    "D29" = "0.0000104"
    "E29" = "  -2.92%  "
    "D30" = "0.999"
    "E30" = "  -0.06%  "
    "D31" = "562.53"
    "E31" = "  -1.47%  "
    "E32" = "  -0.35%  "
    "E33" = "  -1.29%  "
    "E34" = "  -0.15%  "
    "D35" = "1.66"
    "E35" = "  +4.24%  "
    "E36" = "  -1.24%  "
    "E37" = "  -0.03%  "
    "D38" = "160.36"
    "E38" = "  -0.49%  "
    "D39" = "19.65"
    "E39" = "  +1.47%  "
    "E40" = "  -1.28%  "
    "E41" = "  -0.98%  "
    "D42" = "5.34"
    "E42" = "  -0.66%  "
    "D43" = "2.64"
    "E43" = "  -0.60%  "
    "D44" = "0.0₆0322"
    "E44" = "  -4.98%  "
    "D46" = "158.32"
    "E46" = "  +1.11%  "
    "D47" = "3.84"
    "E47" = "  +1.94%  "
    "E48" = "  +0.09%  "
    "D49" = "1.70"
    "E49" = "  -0.92%  "
    "D50" = "0.0777"
    "E50" = "  -1.37%  "
}

foreach ($ref in $cellUpdates.Keys) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $cellUpdates[$ref]
    $rng.Style = "Normal"
}

